$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.176.31'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.058.63'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.86'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.666'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.97'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.76%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0787'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.29%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.22'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.916'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +13.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.360.87'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.75'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.062.25'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.71'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +12.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.210.82'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.93'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0902'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.40'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.50'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.67'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.65%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.50'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.23'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.16'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +9.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.18'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.67%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0626'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.67%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.66'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.79%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.51%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.31'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.77'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.31'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +15.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.09'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +8.15%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -10.58%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.17'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.20%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.56'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.88'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.44'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.278.69'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.00%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.251.19'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.58'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.59%  '
